# Timing issue fix - keywords, updated tc1,2 in ubc01
#
# The "CasesTab" row's query (cell B2) had a stray OPTIONAL MATCH on
# (co:cohort) plus a trailing `Cohort` column that isn't produced by the
# other two queries on this sheet (SamplesTab / FilesTab). Drop the extra
# blank line after the first MATCH and drop the trailing cohort column so
# the case query lines up with the rest of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$caseQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
             "MATCH (c)<--(diag:diagnosis)`n" +
             " MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)`n" +
             "`tWHERE s.clinical_study_designation IN ['UBC01'] and demo.neutered_indicator in [ 'No'] OPTIONAL MATCH (samp:sample)-->(c)`n" +
             "OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" +
             "WITH DISTINCT c, s, demo, diag, co`n" +
             "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" +
             "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" +
             "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" +
             "        coalesce(demo.breed, '') AS Breed ,`n" +
             "        coalesce(diag.disease_term, '') AS Diagnosis ,`n" +
             "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" +
             "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" +
             "        coalesce(demo.sex, '') AS Sex ,`n" +
             "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" +
             "        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" +
             "        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $caseQuery

# The query text got shorter (one fewer wrapped line), so the row now
# needs less height to display fully - matches SamplesTab/FilesTab rows.
$ws.Rows.Item(2).RowHeight = 290

# Move the active selection from the FilesTab row up to the CasesTab row.
$ws.Range("B2").Select()
